# Fruta / hortaliza, semanal
# Refresh the weekly price records (rows 2-22) with the latest batch of
# market observations: dates, quality, volume and price fields are updated
# in place to match the newest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44488
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = 800
$ws.Range("N2").Value = '$/kilo (volumen en unidades)'
$ws.Range("O2").Value = 'Perú'
$ws.Range("P2").Value = 800
$ws.Range("D3").Value = 44305
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2500
$ws.Range("O3").Value = 'Perú'
$ws.Range("P3").Value = 2500
$ws.Range("D4").Value = 44194
$ws.Range("I4").Value = 'Extra'
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 3500
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3500
$ws.Range("N4").Value = '$/unidad'
$ws.Range("O4").Value = 'Región de O''Higgins'
$ws.Range("P4").Value = 3500
$ws.Range("D5").Value = 44194
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("D6").Value = 44167
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = 5000
$ws.Range("P6").Value = 5000
$ws.Range("D7").Value = 44167
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("J7").Value = 560
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("D8").Value = 44167
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("J8").Value = 450
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 2000
$ws.Range("P8").Value = 2000
$ws.Range("D9").Value = 44483
$ws.Range("J9").Value = 120
$ws.Range("D10").Value = 44217
$ws.Range("I10").Value = 'Extra'
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2500
$ws.Range("P10").Value = 2500
$ws.Range("D11").Value = 44217
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 280
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 2000
$ws.Range("P11").Value = 2000
$ws.Range("D12").Value = 44510
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 800
$ws.Range("L12").Value = 800
$ws.Range("M12").Value = 800
$ws.Range("N12").Value = '$/kilo (volumen en unidades)'
$ws.Range("O12").Value = 'Perú'
$ws.Range("P12").Value = 800
$ws.Range("D13").Value = 44504
$ws.Range("J13").Value = 200
$ws.Range("D14").Value = 44495
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 800
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = 800
$ws.Range("N14").Value = '$/kilo (volumen en unidades)'
$ws.Range("P14").Value = 800
$ws.Range("D15").Value = 44312
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 180
$ws.Range("O15").Value = 'Perú'
$ws.Range("D16").Value = 44491
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 800
$ws.Range("N16").Value = '$/kilo (volumen en unidades)'
$ws.Range("O16").Value = 'Perú'
$ws.Range("P16").Value = 800
$ws.Range("D17").Value = 44223
$ws.Range("H17").Value = 'Americana O Klondike'
$ws.Range("I17").Value = 'Extra'
$ws.Range("J17").Value = 340
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2500
$ws.Range("N17").Value = '$/unidad'
$ws.Range("O17").Value = 'Región de O''Higgins'
$ws.Range("P17").Value = 2500
$ws.Range("D18").Value = 44223
$ws.Range("H18").Value = 'Americana O Klondike'
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 2000
$ws.Range("N18").Value = '$/unidad'
$ws.Range("O18").Value = 'Región de O''Higgins'
$ws.Range("P18").Value = 2000
$ws.Range("D19").Value = 44223
$ws.Range("H19").Value = 'Americana O Klondike'
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = 1500
$ws.Range("N19").Value = '$/unidad'
$ws.Range("O19").Value = 'Región de O''Higgins'
$ws.Range("P19").Value = 1500
$ws.Range("D20").Value = 44223
$ws.Range("H20").Value = 'Americana O Klondike'
$ws.Range("I20").Value = 'Tercera'
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = 1000
$ws.Range("O20").Value = 'Región de O''Higgins'
$ws.Range("P20").Value = 1000
$ws.Range("D21").Value = 44477
$ws.Range("J21").Value = 80
$ws.Range("D22").Value = 44497
$ws.Range("J22").Value = 250
